$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.404.52'
$ws.Range("E2").Value = '''  -3.51%  '
$ws.Range("D3").Value = '''1.850.75'
$ws.Range("E3").Value = '''  -5.33%  '
$ws.Range("E4").Value = '''  -0.78%  '
$ws.Range("D5").Value = '''321.46'
$ws.Range("E5").Value = '''  +0.17%  '
$ws.Range("E6").Value = '''  -0.69%  '
$ws.Range("D7").Value = '''0.4459'
$ws.Range("E7").Value = '''  -6.24%  '
$ws.Range("D8").Value = '''0.3833'
$ws.Range("E8").Value = '''  -5.42%  '
$ws.Range("D9").Value = '''49.21'
$ws.Range("E9").Value = '''  -8.01%  '
$ws.Range("D10").Value = '''0.07828'
$ws.Range("E10").Value = '''  -7.30%  '
$ws.Range("E11").Value = '''  -3.89%  '
$ws.Range("D12").Value = '''21.42'
$ws.Range("E12").Value = '''  -3.28%  '
$ws.Range("D13").Value = '''1.839.37'
$ws.Range("E13").Value = '''  -5.99%  '
$ws.Range("D14").Value = '''5.838'
$ws.Range("E14").Value = '''  -4.96%  '
$ws.Range("D15").Value = '''7.099'
$ws.Range("E15").Value = '''  -6.51%  '
$ws.Range("D16").Value = '''1.003'
$ws.Range("E16").Value = '''  -0.74%  '
$ws.Range("D17").Value = '''0.00001026'
$ws.Range("E17").Value = '''  -4.13%  '
$ws.Range("D18").Value = '''85.22'
$ws.Range("E18").Value = '''  -4.65%  '
$ws.Range("D19").Value = '''0.06494'
$ws.Range("E19").Value = '''  -1.60%  '
$ws.Range("D20").Value = '''16.93'
$ws.Range("E20").Value = '''  -8.96%  '
$ws.Range("D21").Value = '''1.002'
$ws.Range("D22").Value = '''5.466'
$ws.Range("E22").Value = '''  -5.90%  '
$ws.Range("D23").Value = '''27.399.23'
$ws.Range("E23").Value = '''  -3.63%  '
$ws.Range("E24").Value = '''  -6.70%  '
$ws.Range("D25").Value = '''2.259'
$ws.Range("E25").Value = '''  -1.12%  '
$ws.Range("D26").Value = '''2.067.16'
$ws.Range("E26").Value = '''  -5.62%  '
$ws.Range("D27").Value = '''151.27'
$ws.Range("E27").Value = '''  -2.25%  '
$ws.Range("D28").Value = '''19.31'
$ws.Range("E28").Value = '''  -4.33%  '
$ws.Range("D29").Value = '''2.044'
$ws.Range("E29").Value = '''  -4.87%  '
$ws.Range("D30").Value = '''5.467'
$ws.Range("E30").Value = '''  -7.52%  '
$ws.Range("D31").Value = '''120.03'
$ws.Range("E31").Value = '''  -2.79%  '
$ws.Range("B32").Value = '''Stellar'
$ws.Range("C32").Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '''0.09309'
$ws.Range("E32").Value = '''  -2.84%  '
$ws.Range("B33").Value = '''ARBITRUM'
$ws.Range("C33").Value = '''https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").Value = '''1.472'
$ws.Range("E33").Value = '''  +2.37%  '
$ws.Range("D34").Value = '''0.9260'
$ws.Range("E34").Value = '''  -5.31%  '
$ws.Range("D35").Value = '''3.625'
$ws.Range("E35").Value = '''  -0.93%  '
$ws.Range("D36").Value = '''5.211'
$ws.Range("E36").Value = '''  -6.52%  '
$ws.Range("D37").Value = '''0.02215'
$ws.Range("E37").Value = '''  -4.88%  '
$ws.Range("D38").Value = '''0.05940'
$ws.Range("E38").Value = '''  -4.32%  '
$ws.Range("D39").Value = '''1.205'
$ws.Range("E39").Value = '''  -3.61%  '
$ws.Range("D40").Value = '''8.314'
$ws.Range("E40").Value = '''  -5.77%  '
$ws.Range("E41").Value = '''  -0.78%  '
$ws.Range("D42").Value = '''0.5908'
$ws.Range("E42").Value = '''  -4.65%  '
$ws.Range("E43").Value = '''  -3.56%  '
$ws.Range("D44").Value = '''10.24'
$ws.Range("E44").Value = '''  -7.70%  '
$ws.Range("D45").Value = '''1.252'
$ws.Range("E45").Value = '''  -6.57%  '
$ws.Range("D46").Value = '''0.5661'
$ws.Range("E46").Value = '''  -4.86%  '
$ws.Range("D47").Value = '''12.11'
$ws.Range("E47").Value = '''  -6.60%  '
$ws.Range("D48").Value = '''3.358'
$ws.Range("E48").Value = '''  -1.04%  '
$ws.Range("D49").Value = '''1.916'
$ws.Range("E49").Value = '''  -6.70%  '
$ws.Range("D50").Value = '''0.06854'
$ws.Range("E50").Value = '''  +0.60%  '
$ws.Range("D51").Value = '''107.97'
$ws.Range("E51").Value = '''  -2.23%  '
